# Update "想去人数" (column F) counts across the 4 sheets of the workbook
# to match the newly-published snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 64
$ws.Range("F6").Value = 915
$ws.Range("F7").Value = 465
$ws.Range("F9").Value = 2194
$ws.Range("F10").Value = 630
$ws.Range("F11").Value = 291
$ws.Range("F12").Value = 120
$ws.Range("F13").Value = 1090
$ws.Range("F14").Value = 183
$ws.Range("F15").Value = 2213
$ws.Range("F16").Value = 674
$ws.Range("F17").Value = 13303
$ws.Range("F18").Value = 1269
$ws.Range("F19").Value = 33
$ws.Range("F20").Value = 561
$ws.Range("F23").Value = 143
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 41
$ws.Range("F26").Value = 271
$ws.Range("F29").Value = 24

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 19
$ws.Range("F11").Value = 83
$ws.Range("F17").Value = 12
$ws.Range("F22").Value = 5

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5711
$ws.Range("F3").Value = 484

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 64
$ws.Range("F3").Value = 484
$ws.Range("F8").Value = 915
$ws.Range("F10").Value = 465
$ws.Range("F11").Value = 19
$ws.Range("F12").Value = 2194
$ws.Range("F13").Value = 630
$ws.Range("F14").Value = 291
$ws.Range("F16").Value = 120
$ws.Range("F18").Value = 1090
$ws.Range("F20").Value = 183
$ws.Range("F23").Value = 2213
$ws.Range("F24").Value = 674
$ws.Range("F25").Value = 83
$ws.Range("F27").Value = 1269
$ws.Range("F28").Value = 33
$ws.Range("F29").Value = 561
$ws.Range("F32").Value = 143
$ws.Range("F33").Value = 3
$ws.Range("F35").Value = 41
$ws.Range("F38").Value = 271
$ws.Range("F45").Value = 5
$ws.Range("F49").Value = 24
